$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 4 to 2 -- new sample data
$ws.Range("A2:A11").Value = 2

# Update the active selection on the sheet to reflect where the user left off
$ws.Range("B13").Select()
